# Weekly update: insert 3 new price rows above the existing data block.
# This pushes the previous rows 43:47 down to 46:50, and the 3 freshly
# inserted rows (43:45) are then filled in with this week's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 43 (shifts old rows 43:47 down to 46:50)
$ws.Range("A43:T45").Insert()

# ---- New row 43 ----
$ws.Cells.Item(43, 1).Value = 3
$ws.Cells.Item(43, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 44468
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100107
$ws.Cells.Item(43, 8).Value = "Otros"
$ws.Cells.Item(43, 9).Value = 100107002
$ws.Cells.Item(43, 10).Value = "Chirimoya"
$ws.Cells.Item(43, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(43, 12).Value = "Especial"
$ws.Cells.Item(43, 13).Value = 45
$ws.Cells.Item(43, 14).Value = 27000
$ws.Cells.Item(43, 15).Value = 27000
$ws.Cells.Item(43, 16).Value = 27000
$ws.Cells.Item(43, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(43, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 19).Value = 2700
$ws.Cells.Item(43, 20).Value = 10

# ---- New row 44 ----
$ws.Cells.Item(44, 1).Value = 3
$ws.Cells.Item(44, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(44, 3).Value = "Coquimbo"
$ws.Cells.Item(44, 4).Value = 44468
$ws.Cells.Item(44, 5).Value = 5
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100107
$ws.Cells.Item(44, 8).Value = "Otros"
$ws.Cells.Item(44, 9).Value = 100107002
$ws.Cells.Item(44, 10).Value = "Chirimoya"
$ws.Cells.Item(44, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 48
$ws.Cells.Item(44, 14).Value = 25000
$ws.Cells.Item(44, 15).Value = 25000
$ws.Cells.Item(44, 16).Value = 25000
$ws.Cells.Item(44, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(44, 19).Value = 2500
$ws.Cells.Item(44, 20).Value = 10

# ---- New row 45 ----
$ws.Cells.Item(45, 1).Value = 3
$ws.Cells.Item(45, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(45, 3).Value = "Coquimbo"
$ws.Cells.Item(45, 4).Value = 44468
$ws.Cells.Item(45, 5).Value = 5
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100107
$ws.Cells.Item(45, 8).Value = "Otros"
$ws.Cells.Item(45, 9).Value = 100107002
$ws.Cells.Item(45, 10).Value = "Chirimoya"
$ws.Cells.Item(45, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(45, 12).Value = "Segunda"
$ws.Cells.Item(45, 13).Value = 40
$ws.Cells.Item(45, 14).Value = 22000
$ws.Cells.Item(45, 15).Value = 22000
$ws.Cells.Item(45, 16).Value = 22000
$ws.Cells.Item(45, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(45, 19).Value = 2200
$ws.Cells.Item(45, 20).Value = 10
